# Ampliación 5 (Ejercicio 1): mark "5. Mejorar la función buyTicket del
# contrato inteligente (0,4)" as done (highlight red -> green), and bump
# the running point total for Ejercicio 1 from "1,3 puntos" to
# "1,7 puntos".

$d = $word.ActiveDocument

# --- Change 1: red -> green highlight on the "buyTicket" bullet ---------
# Confirm the target bullet is present, then flip every run's highlight
# from red to green. Going through the paragraph Range's Font object (as
# opposed to the Range itself) also updates the trailing paragraph-mark's
# own highlight, matching the <w:pPr><w:rPr> entry in the diff.
$found = $d.Content.Find.Execute(
    "5. Mejorar la función buyTicket del contrato inteligente (0,4)",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*Mejorar la función buyTicket del contrato inteligente*") {
            $p.Range.Font.HighlightColorIndex = 4   # wdBrightGreen
        }
    }
}

# --- Change 2: "Total: 1,3 puntos" -> "Total: 1,7 puntos" ---------------
# Scope the Find to the one paragraph that has this exact running total so
# the other "Total: 3,4 puntos" line (Ejercicio 2) and the "3" in the
# table-of-contents page reference are left untouched.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Total:*1,3 puntos*") {
        $p.Range.Find.Execute("3", $false, $false, $false, $false, $false,
            $true, 1, $false, "7", 2)
    }
}
